# Daily attendance processing - reorder the "Recorded By" audit trail
# (column G) so entries read most-recent-first: reverse the
# comma-separated list of names/emails for every row where "System"
# appears alongside another recorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,11,12,13,14,15,29,30,32,33,38,39,40,41,42,56,57,58,59,60,65,66,67,68,69,84,85,89,93,110,111,115,119,136,137,141,145)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Value2
    if ($current -ne $null -and $current -like "*,*") {
        $parts = $current -split ",\s*"
        $n = $parts.Count
        $reversedParts = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversedParts)
    }
}
